# Freezing Fajr and midnight (and the knock-on shifts in A.lasku / Maghrib)
# for the Helsinki "4" adhan sheet: a handful of prayer-time cells are
# retimed, and the placeholder "*" entries in the Fajr column (rows 28-31,
# the "white nights" rows with no astronomical Fajr) are replaced with a
# concrete frozen time instead of the text marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Time($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.NumberFormat = "h:mm"
}

# Row 2
Set-Time "E2" 0.8354166666666667
Set-Time "F2" 0.8534722222222222

# Row 3
Set-Time "B3" 0.18194444444444444
Set-Time "D3" 0.5576388888888889
Set-Time "E3" 0.8368055555555556
Set-Time "F3" 0.8555555555555555

# Row 4
Set-Time "E4" 0.8388888888888889
Set-Time "F4" 0.8569444444444444

# Row 5
Set-Time "E5" 0.8402777777777778
Set-Time "F5" 0.8590277777777777

# Row 6
Set-Time "B6" 0.17291666666666666
Set-Time "D6" 0.5569444444444445
Set-Time "E6" 0.8423611111111111
Set-Time "F6" 0.8604166666666667

# Row 7
Set-Time "D7" 0.5569444444444445
Set-Time "E7" 0.84375
Set-Time "F7" 0.8625

# Row 8
Set-Time "E8" 0.8458333333333333
Set-Time "F8" 0.8645833333333334

# Row 9
Set-Time "E9" 0.8472222222222222
Set-Time "F9" 0.8659722222222223

# Row 10
Set-Time "D10" 0.55625
Set-Time "E10" 0.8493055555555555
Set-Time "F10" 0.8680555555555556

# Row 11
Set-Time "E11" 0.8506944444444444
Set-Time "F11" 0.8694444444444445

# Rows 28-31: replace the "*" placeholder in the Fajr (B) column with the
# frozen time used on the other "no true Fajr" rows.
Set-Time "B28" 0.07361111111111111
Set-Time "B29" 0.07361111111111111
Set-Time "B30" 0.07361111111111111
Set-Time "B31" 0.07361111111111111
